$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASK SUMMARY SHEET")

# Row 3: first task entry (Stage, Task)
$ws.Range("A3").Value = "Project Planning"
$ws.Range("B3").Value = "Prepare Project Plan for client, work on introduction"

# Header row: Name and Week
$ws.Range("C1").Value = "Richard Dobson"
$ws.Range("E1").Value = 2

# Row 4: second task entry (Task)
$ws.Range("B4").Value = "Write product specification for metadata extraction program"

# Remaining row 3/4 data
$ws.Range("A4").Value = "Project Planning"
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 10
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 10

# Update selection to match final state
$ws.Range("B4").Select()
